$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOIL_parameters")

# Row 3: B3 3 -> 1, C3 0 -> 10 with new numeric style (center aligned integer format)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 10
$ws.Range("C3").NumberFormat = "0"
$ws.Range("C3").HorizontalAlignment = -4108  # xlCenter

# Row 4: clear contents (Clay row removed), keep formatting
$ws.Range("A4:G4").ClearContents()

# Update selection on the sheet to C7
$ws.Range("C7").Select()
